$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes existing rows 8-20 down to 9-21,
# keeping their label+data pairing intact)
$ws.Rows.Item(8).Insert()

$c = $ws.Cells.Item(8, 1)

# Match the bold/bordered/centered formatting used by the other A-column
# label cells (mirrors cell style "s=1" from the rest of the sheet)
$c.Font.Bold = $true
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160

# Fill the new row 8 with the new "old IBES data" comparison row
$c.Value = "ibes_2|fwdepsqcut|ibes_industry_all x -exclude_stock"
$ws.Cells.Item(8, 2).Value = 0.006723070191127427
$ws.Cells.Item(8, 3).Value = 0.006502442371883452
$ws.Cells.Item(8, 4).Value = 0.0001075952021389975
$ws.Cells.Item(8, 5).Value = 0.000110182239419974
$ws.Cells.Item(8, 6).Value = 0.224562807151098
$ws.Cells.Item(8, 7).Value = 0.2059180638254218
$ws.Cells.Item(8, 8).Value = 0.301426345673243
$ws.Cells.Item(8, 9).Value = 3343
